# The workbook already holds the "before" state of a weekly price sheet for
# "Hortaliza, Macroferia Regional de Talca - Zapallo italiano".
# A new weekly record needs to be inserted as row 166 (pushing the existing
# rows 166-274 down to 167-275), matching the author's commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 166; everything below (old 166..274) shifts
# down to 167..275, which also grows the used range to A1:R275.
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(166, 1).Value = 5
$ws.Cells.Item(166, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(166, 3).Value = 'Maule'
$ws.Cells.Item(166, 4).Value = 44582
$ws.Cells.Item(166, 5).Value = 7
$ws.Cells.Item(166, 6).Value = 100112032
$ws.Cells.Item(166, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(166, 8).Value = 'Sin especificar'
$ws.Cells.Item(166, 9).Value = 'Primera'
$ws.Cells.Item(166, 10).Value = 850
$ws.Cells.Item(166, 11).Value = 8000
$ws.Cells.Item(166, 12).Value = 8000
$ws.Cells.Item(166, 13).Value = 8000
$ws.Cells.Item(166, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(166, 15).Value = 'Región del Maule'
$ws.Cells.Item(166, 16).Value = 133
$ws.Cells.Item(166, 17).Value = 60
$ws.Cells.Item(166, 18).Value = 'Hortaliza'

# Match the date-serial number format used by the rest of column D.
$ws.Cells.Item(166, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
